# Weekly update: insert a new price record for "Terminal La Palmera de La
# Serena - Papa" right before the former row 571, shifting the existing
# rows 571-612 down to 572-613 (matches the diff: a brand-new data row
# appears, and every later row keeps its old contents one row further
# down).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 571..612 down to 572..613, leaving a blank row 571 to fill in.
$ws.Rows.Item(571).Insert()

$ws.Cells.Item(571, 1).Value  = 8
$ws.Cells.Item(571, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(571, 3).Value  = "Coquimbo"
$ws.Cells.Item(571, 4).Value  = 45013
$ws.Cells.Item(571, 5).Value  = 4
$ws.Cells.Item(571, 6).Value  = 100114001
$ws.Cells.Item(571, 7).Value  = "Papa"
$ws.Cells.Item(571, 8).Value  = "Asterix"
$ws.Cells.Item(571, 9).Value  = "1a (cosecha)"
$ws.Cells.Item(571, 10).Value = 1800
$ws.Cells.Item(571, 11).Value = 11000
$ws.Cells.Item(571, 12).Value = 12000
$ws.Cells.Item(571, 13).Value = 11500
$ws.Cells.Item(571, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(571, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(571, 16).Value = 460
$ws.Cells.Item(571, 17).Value = 25
$ws.Cells.Item(571, 18).Value = "Hortaliza"
